$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.091.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.426.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.99%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.51%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.679"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.52%  "

$ws.Range("E12").Value = "  -1.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.963.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.14%  "

$ws.Range("E14").Value = "  -3.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.413.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.101.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.73%  "

$ws.Range("E18").Value = "  -2.83%  "

$ws.Range("E19").Value = "  +0.66%  "

$ws.Range("E20").Value = "  -2.95%  "

$ws.Range("E21").Value = "  -4.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "84.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "315.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.83%  "

$ws.Range("E24").Value = "  -3.19%  "

$ws.Range("E25").Value = "  -2.45%  "

$ws.Range("E26").Value = "  +9.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.20%  "

$ws.Range("E31").Value = "  -2.03%  "

$ws.Range("E32").Value = "  -4.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "42.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.64%  "

$ws.Range("E34").Value = "  -0.15%  "

$ws.Range("E35").Value = "  -4.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0486"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("E39").Value = "  -4.36%  "

$ws.Range("E40").Value = "  -1.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("E44").Value = "  +0.45%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.80%  "

$ws.Range("E47").Value = "  -2.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.131.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.35%  "

$ws.Range("E50").Value = "  -3.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.97%  "
